# Memphis roster: re-sort a handful of player rows (totals model v1 reorder)
# and backfill Luke Kennard's jersey number.
#
# Columns B..K hold: No., Player, Pos, Ht, Wt, Birth Date, Unnamed:6, Exp,
# College, bbref url. Column A (row index / No. in roster order) is left
# untouched - only each player's own data block (B:K) moves to its new row.
#
# Copy + PasteSpecial (rather than a plain .Value assignment) is used so the
# "Exp" column - which stores digit-only experience values ("1".."9") as
# TEXT - keeps its original text type instead of Excel's normal auto-convert
# of numeric-looking strings into real numbers. PasteSpecial here carries
# only the value/type of the source cell, leaving each destination cell's
# own formatting (borders, hyperlink style, etc.) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch row far outside the used range, used as swap/rotation scratch
# space; cleared at the end of each step so it never lingers in the file.
$scratch = "B200:K200"

function Move-RosterRow($fromRef, $toRef) {
    $ws.Range($fromRef).Copy()
    $ws.Range($toRef).PasteSpecial()
}

# --- Swap row 5 <-> row 6 (Brandon Clarke <-> John Konchar) ---
Move-RosterRow "B5:K5" $scratch
Move-RosterRow "B6:K6" "B5:K5"
Move-RosterRow $scratch "B6:K6"
$ws.Range($scratch).ClearContents()

# --- 4-cycle among rows 9, 10, 11, 12 ---
# after: row9<-row10, row10<-row11, row11<-row12, row12<-row9
Move-RosterRow "B9:K9" $scratch
Move-RosterRow "B10:K10" "B9:K9"
Move-RosterRow "B11:K11" "B10:K10"
Move-RosterRow "B12:K12" "B11:K11"
Move-RosterRow $scratch "B12:K12"
$ws.Range($scratch).ClearContents()

# --- Swap row 13 <-> row 14 (Jake LaRavia <-> Ziaire Williams) ---
Move-RosterRow "B13:K13" $scratch
Move-RosterRow "B14:K14" "B13:K13"
Move-RosterRow $scratch "B14:K14"
$ws.Range($scratch).ClearContents()

# --- Swap row 16 <-> row 17 (Vince Williams Jr. (TW) <-> Kenneth Lofton Jr. (TW)) ---
Move-RosterRow "B16:K16" $scratch
Move-RosterRow "B17:K17" "B16:K16"
Move-RosterRow $scratch "B17:K17"
$ws.Range($scratch).ClearContents()

# --- Backfill Luke Kennard's jersey number (was blank) ---
$ws.Range("B18").Value = 10
